$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 22 - Banglades
$ws.Range("B22").Value = 81523
$ws.Range("C22").Value = 3471
$ws.Range("D22").Value = 17249
$ws.Range("E22").Value = 63179
$ws.Range("G22").Value = 46
$ws.Range("H22").Value = 1095

# Row 39 - Polonia
$ws.Range("B39").Value = 28577
$ws.Range("C39").Value = 376
$ws.Range("D39").Value = 13805
$ws.Range("E39").Value = 13550
$ws.Range("G39").Value = 7
$ws.Range("H39").Value = 1222

# Row 103 - Sri Lanka
$ws.Range("D103").Value = 1196
$ws.Range("E103").Value = 670

# Row 111 - Eslovenia
$ws.Range("B111").Value = 1490
$ws.Range("C111").Value = 2
$ws.Range("E111").Value = 22

# Row 160 - Birmania
$ws.Range("B160").Value = 261
$ws.Range("C160").Value = 1
$ws.Range("E160").Value = 90
